$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 543.3125
$ws.Range("I15").Value = 543.3125
$ws.Range("K15").Value = 1629.9375
$ws.Range("M15").Value = -1460.9375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 50005196
$ws.Range("I64").Value = 83337496
$ws.Range("K64").Value = 83337496
$ws.Range("M64").Value = -83337248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 50005196
$ws.Range("I67").Value = 83337496
$ws.Range("K67").Value = 83337496
$ws.Range("M67").Value = -83336638

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2536.3809
$ws.Range("I98").Value = 2253.3948
$ws.Range("K98").Value = 2253.3948
$ws.Range("M98").Value = -755.3948

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6374.4165
$ws.Range("I116").Value = 5668.3335
$ws.Range("J116").Value = 6609.778
$ws.Range("K116").Value = 5668.3335
$ws.Range("L116").Value = 6609.778
$ws.Range("M116").Value = -2226.3335
$ws.Range("N116").Value = -13493.778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2536.3809
$ws.Range("I122").Value = 2253.3948
$ws.Range("K122").Value = 6760.1844
$ws.Range("M122").Value = -4310.1844

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9557.794
$ws.Range("I132").Value = 6905.161
$ws.Range("K132").Value = 20715.483
$ws.Range("M132").Value = -18185.483

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 11194.728
$ws.Range("J137").Value = 28399.5
$ws.Range("L137").Value = 85198.5
$ws.Range("N137").Value = -90298.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2251.0203
$ws.Range("I138").Value = 1888.8948
$ws.Range("J138").Value = 2476.6064
$ws.Range("K138").Value = 5666.6844
$ws.Range("L138").Value = 7429.8192
$ws.Range("M138").Value = -526.6844000000001
$ws.Range("N138").Value = -17709.8192

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3023
$ws.Range("I141").Value = 3086.9375
$ws.Range("K141").Value = 9260.8125
$ws.Range("M141").Value = -4080.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H81").Value = 82590.5
$ws.Range("J81").Value = 82590.5
$ws.Range("L81").Value = 82590.5
$ws.Range("N81").Value = -84586.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H84").Value = 82590.5
$ws.Range("J84").Value = 82590.5
$ws.Range("L84").Value = 247771.5
$ws.Range("N84").Value = -257755.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11248.667
$ws.Range("I31").Value = 1063.7333
$ws.Range("J31").Value = 21433.6
$ws.Range("K31").Value = 1063.7333
$ws.Range("L31").Value = 21433.6
$ws.Range("M31").Value = -768.7333000000001
$ws.Range("N31").Value = -22023.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11248.667
$ws.Range("I34").Value = 1063.7333
$ws.Range("J34").Value = 21433.6
$ws.Range("K34").Value = 1063.7333
$ws.Range("L34").Value = 21433.6
$ws.Range("M34").Value = -861.7333000000001
$ws.Range("N34").Value = -21837.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13984.806
$ws.Range("I58").Value = 6684.864
$ws.Range("K58").Value = 6684.864
$ws.Range("M58").Value = -6481.864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2331.6667
$ws.Range("I62").Value = 2331.6667
$ws.Range("K62").Value = 2331.6667
$ws.Range("M62").Value = -1707.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2331.6667
$ws.Range("I65").Value = 2331.6667
$ws.Range("K65").Value = 11658.3335
$ws.Range("M65").Value = -8538.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 16617.666
$ws.Range("I86").Value = 21099
$ws.Range("K86").Value = 21099
$ws.Range("M86").Value = -19976

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 16617.666
$ws.Range("I89").Value = 21099
$ws.Range("K89").Value = 105495
$ws.Range("M89").Value = -99879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 32260582
$ws.Range("I132").Value = 2416.0344
$ws.Range("K132").Value = 7248.1032
$ws.Range("M132").Value = -4718.1032

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 35721572
$ws.Range("I134").Value = 1673.8667
$ws.Range("K134").Value = 5021.6001
$ws.Range("M134").Value = -2486.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13984.806
$ws.Range("I136").Value = 6684.864
$ws.Range("K136").Value = 20054.592
$ws.Range("M136").Value = -17504.592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 18183408
$ws.Range("I129").Value = 2316.3333
$ws.Range("K129").Value = 6948.999899999999
$ws.Range("M129").Value = -1948.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 89995
$ws.Range("J32").Value = 89995
$ws.Range("L32").Value = 89995
$ws.Range("N32").Value = -90587

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17199.4
$ws.Range("I70").Value = 12499.5
$ws.Range("J70").Value = 20332.666
$ws.Range("K70").Value = 12499.5
$ws.Range("L70").Value = 20332.666
$ws.Range("M70").Value = -12229.5
$ws.Range("N70").Value = -20872.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 17199.4
$ws.Range("I73").Value = 12499.5
$ws.Range("J73").Value = 20332.666
$ws.Range("K73").Value = 12499.5
$ws.Range("L73").Value = 20332.666
$ws.Range("M73").Value = -11563.5
$ws.Range("N73").Value = -22204.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3000.0908
$ws.Range("I113").Value = 2478.6
$ws.Range("K113").Value = 2478.6
$ws.Range("M113").Value = -308.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 102247.75
$ws.Range("J135").Value = 102247.75
$ws.Range("L135").Value = 102247.75
$ws.Range("N135").Value = -112387.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1009854.6
$ws.Range("I132").Value = 5013.8696
$ws.Range("J132").Value = 2369345
$ws.Range("K132").Value = 15041.6088
$ws.Range("L132").Value = 7108035
$ws.Range("M132").Value = -12511.6088
$ws.Range("N132").Value = -7113095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 45874.5
$ws.Range("J74").Value = 45874.5
$ws.Range("L74").Value = 45874.5
$ws.Range("N74").Value = -47746.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 45874.5
$ws.Range("J77").Value = 45874.5
$ws.Range("L77").Value = 137623.5
$ws.Range("N77").Value = -146983.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16202.189
$ws.Range("I132").Value = 8258.182000000001
$ws.Range("J132").Value = 27853.4
$ws.Range("K132").Value = 24774.546
$ws.Range("L132").Value = 83560.20000000001
$ws.Range("M132").Value = -22244.546
$ws.Range("N132").Value = -88620.20000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12227.833
$ws.Range("I136").Value = 2005.4642
$ws.Range("J136").Value = 32672.572
$ws.Range("K136").Value = 6016.392599999999
$ws.Range("L136").Value = 98017.716
$ws.Range("M136").Value = -3466.392599999999
$ws.Range("N136").Value = -103117.716
